$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.70949999999999
$ws.Range("A6").Value = -22.4235
$ws.Range("A7").Value = -20.07919999999999
$ws.Range("C7").Value = -12.1903
$ws.Range("A8").Value = -22.31840000000001
$ws.Range("C11").Value = -11.1651
$ws.Range("C12").Value = -10.7955
$ws.Range("D12").Value = -7.119999999999997
$ws.Range("D13").Value = -8.4971
$ws.Range("D14").Value = -7.921399999999997
$ws.Range("C15").Value = -14.54699999999999
$ws.Range("A16").Value = -21.67949999999999
$ws.Range("D16").Value = -9.095400000000005
$ws.Range("D19").Value = -7.9085
$ws.Range("A20").Value = -19.7301
$ws.Range("C20").Value = -11.81360000000001
$ws.Range("D20").Value = -7.054299999999997
$ws.Range("A21").Value = -19.69679999999999
$ws.Range("C21").Value = -11.83380000000001
$ws.Range("C22").Value = -12.48779999999999
$ws.Range("D22").Value = -8.186200000000003
$ws.Range("C23").Value = -12.00480000000001
$ws.Range("A28").Value = -21.75759999999999
$ws.Range("A29").Value = -21.59389999999999
$ws.Range("C29").Value = -11.57310000000001
$ws.Range("A30").Value = -21.6085
$ws.Range("A32").Value = -21.23679999999999
$ws.Range("C34").Value = -11.26350000000001
$ws.Range("D36").Value = -8.248199999999995
$ws.Range("A40").Value = -20.1227
$ws.Range("C42").Value = -12.1103
$ws.Range("C43").Value = -13.6053
$ws.Range("D43").Value = -8.354800000000001
$ws.Range("C44").Value = -13.7311
$ws.Range("C45").Value = -13.74479999999999
$ws.Range("A46").Value = -21.9148
$ws.Range("C46").Value = -12.60370000000001
$ws.Range("D46").Value = -8.688899999999995
$ws.Range("C50").Value = -14.16869999999998
$ws.Range("D50").Value = -7.978400000000001
$ws.Range("A51").Value = -21.5848
$ws.Range("C51").Value = -10.89690000000001
$ws.Range("A52").Value = -22.2895
$ws.Range("A57").Value = -21.8157
$ws.Range("C57").Value = -11.95069999999999
$ws.Range("A59").Value = -22.25960000000001
$ws.Range("A62").Value = -22.3011
$ws.Range("C65").Value = -13.07369999999999
$ws.Range("A66").Value = -21.581
$ws.Range("C66").Value = -11.01430000000001
$ws.Range("C67").Value = -11.3015
$ws.Range("A73").Value = -20.5076
$ws.Range("A74").Value = -22.03869999999998
$ws.Range("D76").Value = -7.804299999999999
$ws.Range("A77").Value = -19.66769999999999
$ws.Range("C79").Value = -11.69570000000001
$ws.Range("C84").Value = -12.86919999999999
$ws.Range("C87").Value = -13.6093
$ws.Range("A92").Value = -21.52020000000001
$ws.Range("C92").Value = -11.5321
$ws.Range("D95").Value = -8.221999999999996
$ws.Range("C97").Value = -11.0643
$ws.Range("D97").Value = -8.260599999999993
$ws.Range("D99").Value = -8.101099999999997
$ws.Range("A100").Value = -22.2338
